$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 769.05884
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 692.125
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 692.125
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -830.125
$ws.Range("H62").Value = 2529.2
$ws.Range("I62").Value = 1468.5714
$ws.Range("J62").Value = 5004
$ws.Range("K62").Value = 1468.5714
$ws.Range("L62").Value = 5004
$ws.Range("M62").Value = -844.5714
$ws.Range("N62").Value = -6252
$ws.Range("H65").Value = 2529.2
$ws.Range("I65").Value = 1468.5714
$ws.Range("J65").Value = 5004
$ws.Range("K65").Value = 7342.857
$ws.Range("L65").Value = 25020
$ws.Range("M65").Value = -4222.857
$ws.Range("N65").Value = -31260
$ws.Range("H98").Value = 2631.1428
$ws.Range("I98").Value = 1372.4
$ws.Range("J98").Value = 5778
$ws.Range("K98").Value = 1372.4
$ws.Range("L98").Value = 5778
$ws.Range("M98").Value = 125.5999999999999
$ws.Range("N98").Value = -8774
$ws.Range("H115").Value = 1610.3125
$ws.Range("I115").Value = 988.75
$ws.Range("J115").Value = 3475
$ws.Range("K115").Value = 2966.25
$ws.Range("L115").Value = 10425
$ws.Range("M115").Value = -1399.25
$ws.Range("N115").Value = -13559
$ws.Range("H122").Value = 2631.1428
$ws.Range("I122").Value = 1372.4
$ws.Range("J122").Value = 5778
$ws.Range("K122").Value = 4117.200000000001
$ws.Range("L122").Value = 17334
$ws.Range("M122").Value = -1667.200000000001
$ws.Range("N122").Value = -22234.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 855.9697
$ws.Range("I94").Value = 780.41174
$ws.Range("J94").Value = 936.25
$ws.Range("K94").Value = 780.41174
$ws.Range("L94").Value = 936.25
$ws.Range("M94").Value = -329.41174
$ws.Range("N94").Value = -1838.25
$ws.Range("H105").Value = 2432.5625
$ws.Range("I105").Value = 2005
$ws.Range("J105").Value = 2493.6428
$ws.Range("K105").Value = 2005
$ws.Range("L105").Value = 2493.6428
$ws.Range("M105").Value = -258
$ws.Range("N105").Value = -5987.6428
$ws.Range("H107").Value = 3629.6843
$ws.Range("I107").Value = 2580.5
$ws.Range("J107").Value = 6567.4
$ws.Range("K107").Value = 2580.5
$ws.Range("L107").Value = 6567.4
$ws.Range("M107").Value = -660.5
$ws.Range("N107").Value = -10407.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2348.8
$ws.Range("I31").Value = 1752.909
$ws.Range("J31").Value = 3987.5
$ws.Range("K31").Value = 1752.909
$ws.Range("L31").Value = 3987.5
$ws.Range("M31").Value = -1457.909
$ws.Range("N31").Value = -4577.5
$ws.Range("H34").Value = 2348.8
$ws.Range("I34").Value = 1752.909
$ws.Range("J34").Value = 3987.5
$ws.Range("K34").Value = 1752.909
$ws.Range("L34").Value = 3987.5
$ws.Range("M34").Value = -1550.909
$ws.Range("N34").Value = -4391.5
$ws.Range("H107").Value = 1928.6428
$ws.Range("I107").Value = 778
$ws.Range("J107").Value = 2388.9
$ws.Range("K107").Value = 778
$ws.Range("L107").Value = 2388.9
$ws.Range("M107").Value = 1142
$ws.Range("N107").Value = -6228.9
$ws.Range("H132").Value = 3288.2334
$ws.Range("I132").Value = 3043.1365
$ws.Range("J132").Value = 3962.25
$ws.Range("K132").Value = 9129.4095
$ws.Range("L132").Value = 11886.75
$ws.Range("M132").Value = -6599.4095
$ws.Range("N132").Value = -16946.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3409.0908
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3409.0908
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10227.2724
$ws.Range("N51").Value = -11147.2724
$ws.Range("M51").ClearContents()
$ws.Range("H70").Value = 6645.8335
$ws.Range("I70").Value = 5990
$ws.Range("J70").Value = 6777
$ws.Range("K70").Value = 17970
$ws.Range("L70").Value = 20331
$ws.Range("M70").Value = -17655
$ws.Range("N70").Value = -20961
$ws.Range("H73").Value = 6645.8335
$ws.Range("I73").Value = 5990
$ws.Range("J73").Value = 6777
$ws.Range("K73").Value = 17970
$ws.Range("L73").Value = 20331
$ws.Range("M73").Value = -16878
$ws.Range("N73").Value = -22515
$ws.Range("H75").Value = 3559.1428
$ws.Range("I75").Value = 604.3333
$ws.Range("J75").Value = 4365
$ws.Range("K75").Value = 1812.9999
$ws.Range("L75").Value = 13095
$ws.Range("M75").Value = -814.9999
$ws.Range("N75").Value = -15091
$ws.Range("H78").Value = 3559.1428
$ws.Range("I78").Value = 604.3333
$ws.Range("J78").Value = 4365
$ws.Range("K78").Value = 5438.9997
$ws.Range("L78").Value = 39285
$ws.Range("M78").Value = -446.9997000000003
$ws.Range("N78").Value = -49269
$ws.Range("H131").Value = 888.74
$ws.Range("I131").Value = 643.6667
$ws.Range("J131").Value = 896.3196
$ws.Range("K131").Value = 1931.0001
$ws.Range("L131").Value = 2688.9588
$ws.Range("M131").Value = 3108.9999
$ws.Range("N131").Value = -12768.9588

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 69900
$ws.Range("J28").Value = 69900
$ws.Range("L28").Value = 69900
$ws.Range("N28").Value = -70284
$ws.Range("H122").Value = 101808.8
$ws.Range("I122").Value = 201956
$ws.Range("J122").Value = 1661.6
$ws.Range("K122").Value = 605868
$ws.Range("L122").Value = 4984.799999999999
$ws.Range("M122").Value = -603418
$ws.Range("N122").Value = -9884.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 26803.8
$ws.Range("I4").Value = 4001.1428
$ws.Range("J4").Value = 80010
$ws.Range("K4").Value = 4001.1428
$ws.Range("L4").Value = 80010
$ws.Range("M4").Value = -3888.1428
$ws.Range("N4").Value = -80236
$ws.Range("H5").Value = 3500
$ws.Range("I5").Value = 3500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3387
$ws.Range("N5").ClearContents()
$ws.Range("H21").Value = 70007
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70007
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70007
$ws.Range("N21").Value = -70355
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 3000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = -2770
$ws.Range("H26").Value = 10377.5
$ws.Range("H28").Value = 26803.8
$ws.Range("I28").Value = 4001.1428
$ws.Range("J28").Value = 80010
$ws.Range("K28").Value = 4001.1428
$ws.Range("L28").Value = 80010
$ws.Range("M28").Value = -3769.1428
$ws.Range("N28").Value = -80474
$ws.Range("H30").Value = 32500
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 32500
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 32500
$ws.Range("N30").Value = -32716
$ws.Range("M30").ClearContents()
$ws.Range("H31").Value = 327.25
$ws.Range("I31").Value = 327.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 327.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -79.25
$ws.Range("N31").ClearContents()
$ws.Range("H33").Value = 8507.5
$ws.Range("I33").Value = 5015
$ws.Range("K33").Value = 5015
$ws.Range("M33").Value = -4725
$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10344
$ws.Range("H37").Value = 26803.8
$ws.Range("I37").Value = 4001.1428
$ws.Range("J37").Value = 80010
$ws.Range("K37").Value = 4001.1428
$ws.Range("L37").Value = 80010
$ws.Range("M37").Value = -3894.1428
$ws.Range("N37").Value = -80224
$ws.Range("H40").Value = 8857
$ws.Range("I40").Value = 8857
$ws.Range("K40").Value = 8857
$ws.Range("M40").Value = -8721
$ws.Range("H93").Value = 19835.621
$ws.Range("I93").Value = 1034.5
$ws.Range("J93").Value = 64274.637
$ws.Range("K93").Value = 1034.5
$ws.Range("L93").Value = 64274.637
$ws.Range("M93").Value = 213.5
$ws.Range("N93").Value = -66770.637
$ws.Range("H132").Value = 3342.276
$ws.Range("I132").Value = 2673.2
$ws.Range("K132").Value = 8019.599999999999
$ws.Range("M132").Value = -5489.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 50182.96
$ws.Range("I122").Value = 2216.9333
$ws.Range("J122").Value = 122132
$ws.Range("K122").Value = 6650.7999
$ws.Range("L122").Value = 366396
$ws.Range("M122").Value = -4200.7999
$ws.Range("N122").Value = -371296
